$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column F (shifts F:T -> G:U)
$ws.Columns("F:F").Insert()

$ws.Range("F1").Value = "UBID"
$ws.Range("F2").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F3").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F4").Value = "VQADDOC3+V8E-VQADDOC3+XBF-VQADDOC3+2EE"
$ws.Range("F5").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F6").Value = "ZIUC82DT+4X5-ZIUC82DT+C4M-ZIUC82DT+YK9"
$ws.Range("F7").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F8").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F9").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
$ws.Range("F10").Value = "Z01TDR2Z+7ES-Z01TDR2Z+HX7-Z01TDR2Z+UAX"
